$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name (through date 2022-11-10 -> 2022-11-11)
$ws.Name = "Through 2022-11-11"

# Update header text in B1 to match new through-date
$ws.Range("B1").Value = "November 2022 (through November 11)"

# Update/add data cells (new carjacking data for 2022-11-11)
$ws.Range("M2").Value = 4
$ws.Range("X2").Value = 3

$ws.Range("B3").Value = 1
$ws.Range("BE3").Value = 7
$ws.Range("CA3").Value = 3

$ws.Range("X5").Value = 6

$ws.Range("B9").Value = 1

$ws.Range("X11").Value = 2

$ws.Range("X14").Value = 8
$ws.Range("BE14").Value = 2

$ws.Range("B16").Value = 1

$ws.Range("M17").Value = 2
$ws.Range("X17").Value = 4

$ws.Range("M18").Value = 1
$ws.Range("BP18").Value = 1

$ws.Range("AI21").Value = 1

$ws.Range("B22").Value = 1

$ws.Range("BE26").Value = 7

$ws.Range("BP28").Value = 1

$ws.Range("M31").Value = 2

$ws.Range("BE33").Value = 3

$ws.Range("BP39").Value = 1

$ws.Range("AI50").Value = 1

$ws.Range("M56").Value = 2

$ws.Range("M64").Value = 1

$ws.Range("X90").Value = 2

$ws.Range("BE97").Value = 1
